$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H15").Value = 97.86
$ws.Range("I15").Value = 97.86
$ws.Range("K15").Value = 293.58
$ws.Range("M15").Value = -124.58
$ws.Range("H98").Value = 5858.25
$ws.Range("I98").Value = 3162.375
$ws.Range("K98").Value = 3162.375
$ws.Range("M98").Value = -1664.375
$ws.Range("H112").Value = 1755.2759
$ws.Range("J112").Value = 1872.12
$ws.Range("L112").Value = 5616.36
$ws.Range("N112").Value = -7832.36
$ws.Range("H113").Value = 10131.667
$ws.Range("I113").Value = 916
$ws.Range("J113").Value = 16714.285
$ws.Range("K113").Value = 916
$ws.Range("L113").Value = 16714.285
$ws.Range("M113").Value = 2338
$ws.Range("N113").Value = -23222.285
$ws.Range("H122").Value = 5858.25
$ws.Range("I122").Value = 3162.375
$ws.Range("K122").Value = 9487.125
$ws.Range("M122").Value = -7037.125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 2367
$ws.Range("I127").Value = 991.4286
$ws.Range("J127").Value = 3242.3635
$ws.Range("K127").Value = 2974.2858
$ws.Range("L127").Value = 9727.0905
$ws.Range("M127").Value = 1985.7142
$ws.Range("N127").Value = -19647.0905
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H133").Value = 59790
$ws.Range("J133").Value = 59790
$ws.Range("L133").Value = 59790
$ws.Range("N133").Value = -69910
$ws.Range("H134").Value = 64766.668
$ws.Range("J134").Value = 64766.668
$ws.Range("L134").Value = 64766.668
$ws.Range("N134").Value = -74906.66800000001
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 4314.864
$ws.Range("I138").Value = 2382.4
$ws.Range("J138").Value = 4659.9463
$ws.Range("K138").Value = 7147.200000000001
$ws.Range("L138").Value = 13979.8389
$ws.Range("M138").Value = -2007.200000000001
$ws.Range("N138").Value = -24259.8389
$ws.Range("H139").Value = 38703.934
$ws.Range("J139").Value = 38703.934
$ws.Range("L139").Value = 38703.934
$ws.Range("N139").Value = -48983.934
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1220.5
$ws.Range("I2").Value = 1194.375
$ws.Range("J2").Value = 1325
$ws.Range("K2").Value = 1194.375
$ws.Range("L2").Value = 1325
$ws.Range("M2").Value = -1081.375
$ws.Range("N2").Value = -1551
$ws.Range("H32").Value = 17182.334
$ws.Range("I32").Value = 13683.48
$ws.Range("K32").Value = 13683.48
$ws.Range("M32").Value = -13396.48
$ws.Range("H102").Value = 1940.8
$ws.Range("J102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("N102").Value = -5244
$ws.Range("H111").Value = 34222
$ws.Range("J111").Value = 34222
$ws.Range("L111").Value = 34222
$ws.Range("N111").Value = -42402
$ws.Range("H116").Value = 1220.5
$ws.Range("I116").Value = 1194.375
$ws.Range("J116").Value = 1325
$ws.Range("K116").Value = 1194.375
$ws.Range("L116").Value = 1325
$ws.Range("M116").Value = 1099.625
$ws.Range("N116").Value = -5913
$ws.Range("H122").Value = 4264.25
$ws.Range("I122").Value = 1962.8
$ws.Range("J122").Value = 8100
$ws.Range("K122").Value = 5888.4
$ws.Range("L122").Value = 24300
$ws.Range("M122").Value = -3438.4
$ws.Range("N122").Value = -29200
$ws.Range("H132").Value = 3550.25
$ws.Range("I132").Value = 1822.4286
$ws.Range("K132").Value = 5467.2858
$ws.Range("M132").Value = -2937.2858
$ws.Range("H137").Value = 44267.8
$ws.Range("J137").Value = 44267.8
$ws.Range("L137").Value = 44267.8
$ws.Range("N137").Value = -54467.8
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1220.5
$ws.Range("I3").Value = 1194.375
$ws.Range("J3").Value = 1325
$ws.Range("K3").Value = 1194.375
$ws.Range("L3").Value = 1325
$ws.Range("M3").Value = -1080.375
$ws.Range("N3").Value = -1553
$ws.Range("H86").Value = 2900
$ws.Range("I86").Value = 2060
$ws.Range("J86").Value = 3425
$ws.Range("K86").Value = 2060
$ws.Range("L86").Value = 3425
$ws.Range("M86").Value = -937
$ws.Range("N86").Value = -5671
$ws.Range("H89").Value = 2900
$ws.Range("I89").Value = 2060
$ws.Range("J89").Value = 3425
$ws.Range("K89").Value = 10300
$ws.Range("L89").Value = 17125
$ws.Range("M89").Value = -4684
$ws.Range("N89").Value = -28357
$ws.Range("H134").Value = 3830.9167
$ws.Range("I134").Value = 1956
$ws.Range("J134").Value = 8092.091
$ws.Range("K134").Value = 5868
$ws.Range("L134").Value = 24276.273
$ws.Range("M134").Value = -3333
$ws.Range("N134").Value = -29346.273
$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 2753.5833
$ws.Range("I132").Value = 2164.875
$ws.Range("K132").Value = 6494.625
$ws.Range("M132").Value = -3964.625
$ws.Range("H134").Value = 13196.4
$ws.Range("I134").Value = 12996
$ws.Range("K134").Value = 38988
$ws.Range("M134").Value = -36453
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 2483.2307
$ws.Range("I102").Value = 2172.628
$ws.Range("K102").Value = 2172.628
$ws.Range("M102").Value = -550.6280000000002
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 1416.9
$ws.Range("I113").Value = 1268.1666
$ws.Range("K113").Value = 1268.1666
$ws.Range("M113").Value = 901.8334
$ws.Range("H132").Value = 3419.48
$ws.Range("I132").Value = 1650.25
$ws.Range("J132").Value = 4252.0586
$ws.Range("K132").Value = 4950.75
$ws.Range("L132").Value = 12756.1758
$ws.Range("M132").Value = -2420.75
$ws.Range("N132").Value = -17816.1758
$ws.Range("H141").Value = 57780
$ws.Range("J141").Value = 59336
$ws.Range("L141").Value = 59336
$ws.Range("N141").Value = -69696
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 4223.273
$ws.Range("I136").Value = 1868
$ws.Range("J136").Value = 5853.846
$ws.Range("K136").Value = 5604
$ws.Range("L136").Value = 17561.538
$ws.Range("M136").Value = -3054
$ws.Range("N136").Value = -22661.538
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 6958.5386
$ws.Range("I122").Value = 5091.8335
$ws.Range("K122").Value = 15275.5005
$ws.Range("M122").Value = -12825.5005
$ws.Range("H132").Value = 4447253
$ws.Range("I132").Value = 2787.2642
$ws.Range("J132").Value = 15154375
$ws.Range("K132").Value = 8361.792600000001
$ws.Range("L132").Value = 45463125
$ws.Range("M132").Value = -5831.792600000001
$ws.Range("N132").Value = -45468185
$ws.Range("H136").Value = 5212.1763
$ws.Range("I136").Value = 1564.1818
$ws.Range("J136").Value = 11900.167
$ws.Range("K136").Value = 4692.5454
$ws.Range("L136").Value = 35700.501
$ws.Range("M136").Value = -2142.5454
$ws.Range("N136").Value = -40800.501
